# Team Attendance B8-G1 16-AUG
# Adds the 16-Aug-2023 attendance row (row 10) with PRESENT/ABSENT marks
# for each team member, plus RENUKA's per-person remark comments on the
# ABSENT cells, mirroring the pattern used for the earlier rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New date row: 16-Aug-2023 (serial 45154), formatted like the rows above it.
$ws.Range("A10").Value = 45154
$ws.Range("A10").NumberFormat = "d-mmm"

$ws.Range("B10").Value = "PRESENT"
$ws.Range("C10").Value = "ABSENT"
$ws.Range("D10").Value = "PRESENT"
$ws.Range("E10").Value = "PRESENT"
$ws.Range("F10").Value = "ABSENT"
$ws.Range("G10").Value = "ABSENT"
$ws.Range("H10").Value = "ABSENT"
$ws.Range("I10").Value = "ABSENT"

# RENUKA's remarks explaining the ABSENT marks for this row.
$nl = [char]10
$ws.Range("C10").AddComment("RENUKA:" + $nl + "Not well health issue")
$ws.Range("F10").AddComment("RENUKA:" + $nl + "No Response")
$ws.Range("G10").AddComment("RENUKA:" + $nl + "No Response")
$ws.Range("H10").AddComment("RENUKA:" + $nl + "No Response")
$ws.Range("I10").AddComment("RENUKA:" + $nl + "No Response")

# Leave the selection where the author left it when they saved.
$ws.Range("D19").Select()
